$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Execute" column (C) for rows 11-13: remove Extension Payments tax type from execution
$ws.Range("C11:C13").Value = "DONOTRUN"

# Adjust column C width to fit new content, no longer auto bestFit
$ws.Columns.Item(3).ColumnWidth = 13.83

# Update selection to reflect the edited range
$ws.Range("C11:C13").Select()
